# Update Fbn1-Itgb3 NATMI LR-pair sheet with revised values
# (recomputed per Dr Hou advice: Ligand/Receptor-expressing cells now 3 instead of 1,
# which changes total expression values and derived specificity scores).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 11.98327633333333
$ws.Cells.Item(2, 8).Value = 35.949829
$ws.Cells.Item(2, 9).Value = 0.03345300399843466
$ws.Cells.Item(2, 10).Value = 0.03345300399843466
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 3.778439
$ws.Cells.Item(2, 14).Value = 11.335317
$ws.Cells.Item(2, 15).Value = 0.4252971528324392
$ws.Cells.Item(2, 16).Value = 0.4252971528324392
$ws.Cells.Item(2, 17).Value = 45.27807864564367
$ws.Cells.Item(2, 18).Value = 407.502707810793
$ws.Cells.Item(2, 19).Value = 0.01422746735422646
$ws.Cells.Item(2, 20).Value = 0.01422746735422646

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 11.98327633333333
$ws.Cells.Item(3, 8).Value = 35.949829
$ws.Cells.Item(3, 9).Value = 0.03345300399843466
$ws.Cells.Item(3, 10).Value = 0.03345300399843466
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 4.333403333333333
$ws.Cells.Item(3, 14).Value = 13.00021
$ws.Cells.Item(3, 15).Value = 0.4877633593505858
$ws.Cells.Item(3, 16).Value = 0.4877633593505858
$ws.Cells.Item(3, 17).Value = 51.92836960712111
$ws.Cells.Item(3, 18).Value = 467.35532646409
$ws.Cells.Item(3, 19).Value = 0.01631714961064507
$ws.Cells.Item(3, 20).Value = 0.01631714961064507

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 11.98327633333333
$ws.Cells.Item(4, 8).Value = 35.949829
$ws.Cells.Item(4, 9).Value = 0.03345300399843466
$ws.Cells.Item(4, 10).Value = 0.03345300399843466
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.2909853333333334
$ws.Cells.Item(4, 14).Value = 0.8729560000000001
$ws.Cells.Item(4, 15).Value = 0.03275300561492853
$ws.Cells.Item(4, 16).Value = 0.03275300561492853
$ws.Cells.Item(4, 17).Value = 3.486957658280445
$ws.Cells.Item(4, 18).Value = 31.382618924524
$ws.Cells.Item(4, 19).Value = 0.001095686427796957
$ws.Cells.Item(4, 20).Value = 0.001095686427796957

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 11.98327633333333
$ws.Cells.Item(5, 8).Value = 35.949829
$ws.Cells.Item(5, 9).Value = 0.03345300399843466
$ws.Cells.Item(5, 10).Value = 0.03345300399843466
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.4814053333333333
$ws.Cells.Item(5, 14).Value = 1.444216
$ws.Cells.Item(5, 15).Value = 0.0541864822020464
$ws.Cells.Item(5, 16).Value = 0.05418648220204641
$ws.Cells.Item(5, 17).Value = 5.768813137673778
$ws.Cells.Item(5, 18).Value = 51.919318239064
$ws.Cells.Item(5, 19).Value = 0.001812700605766167
$ws.Cells.Item(5, 20).Value = 0.001812700605766167

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 332.3726806666667
$ws.Cells.Item(6, 8).Value = 997.1180420000001
$ws.Cells.Item(6, 9).Value = 0.9278651602470024
$ws.Cells.Item(6, 10).Value = 0.9278651602470025
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 3.778439
$ws.Cells.Item(6, 14).Value = 11.335317
$ws.Cells.Item(6, 15).Value = 0.4252971528324392
$ws.Cells.Item(6, 16).Value = 0.4252971528324392
$ws.Cells.Item(6, 17).Value = 1255.849899165479
$ws.Cells.Item(6, 18).Value = 11302.64909248932
$ws.Cells.Item(6, 19).Value = 0.3946184108654651
$ws.Cells.Item(6, 20).Value = 0.3946184108654652

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 332.3726806666667
$ws.Cells.Item(7, 8).Value = 997.1180420000001
$ws.Cells.Item(7, 9).Value = 0.9278651602470024
$ws.Cells.Item(7, 10).Value = 0.9278651602470025
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 4.333403333333333
$ws.Cells.Item(7, 14).Value = 13.00021
$ws.Cells.Item(7, 15).Value = 0.4877633593505858
$ws.Cells.Item(7, 16).Value = 0.4877633593505858
$ws.Cells.Item(7, 17).Value = 1440.304882309869
$ws.Cells.Item(7, 18).Value = 12962.74394078882
$ws.Cells.Item(7, 19).Value = 0.4525786275864475
$ws.Cells.Item(7, 20).Value = 0.4525786275864476

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 332.3726806666667
$ws.Cells.Item(8, 8).Value = 997.1180420000001
$ws.Cells.Item(8, 9).Value = 0.9278651602470024
$ws.Cells.Item(8, 10).Value = 0.9278651602470025
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 0.2909853333333334
$ws.Cells.Item(8, 14).Value = 0.8729560000000001
$ws.Cells.Item(8, 15).Value = 0.03275300561492853
$ws.Cells.Item(8, 16).Value = 0.03275300561492853
$ws.Cells.Item(8, 17).Value = 96.71557527468357
$ws.Cells.Item(8, 18).Value = 870.4401774721521
$ws.Cells.Item(8, 19).Value = 0.03039037280346663
$ws.Cells.Item(8, 20).Value = 0.03039037280346664

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 332.3726806666667
$ws.Cells.Item(9, 8).Value = 997.1180420000001
$ws.Cells.Item(9, 9).Value = 0.9278651602470024
$ws.Cells.Item(9, 10).Value = 0.9278651602470025
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.4814053333333333
$ws.Cells.Item(9, 14).Value = 1.444216
$ws.Cells.Item(9, 15).Value = 0.0541864822020464
$ws.Cells.Item(9, 16).Value = 0.05418648220204641
$ws.Cells.Item(9, 17).Value = 160.0059811272302
$ws.Cells.Item(9, 18).Value = 1440.053830145072
$ws.Cells.Item(9, 19).Value = 0.05027774899162313
$ws.Cells.Item(9, 20).Value = 0.05027774899162314

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 0.08615933333333332
$ws.Cells.Item(10, 8).Value = 0.258478
$ws.Cells.Item(10, 9).Value = 0.0002405259164795302
$ws.Cells.Item(10, 10).Value = 0.0002405259164795302
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 3.778439
$ws.Cells.Item(10, 14).Value = 11.335317
$ws.Cells.Item(10, 15).Value = 0.4252971528324392
$ws.Cells.Item(10, 16).Value = 0.4252971528324392
$ws.Cells.Item(10, 17).Value = 0.3255477852806666
$ws.Cells.Item(10, 18).Value = 2.929930067526
$ws.Cells.Item(10, 19).Value = 0.0001022949874611573
$ws.Cells.Item(10, 20).Value = 0.0001022949874611573

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 0.08615933333333332
$ws.Cells.Item(11, 8).Value = 0.258478
$ws.Cells.Item(11, 9).Value = 0.0002405259164795302
$ws.Cells.Item(11, 10).Value = 0.0002405259164795302
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 4.333403333333333
$ws.Cells.Item(11, 14).Value = 13.00021
$ws.Cells.Item(11, 15).Value = 0.4877633593505858
$ws.Cells.Item(11, 16).Value = 0.4877633593505858
$ws.Cells.Item(11, 17).Value = 0.3733631422644444
$ws.Cells.Item(11, 18).Value = 3.36026828038
$ws.Cells.Item(11, 19).Value = 0.0001173197290329341
$ws.Cells.Item(11, 20).Value = 0.0001173197290329341

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 0.08615933333333332
$ws.Cells.Item(12, 8).Value = 0.258478
$ws.Cells.Item(12, 9).Value = 0.0002405259164795302
$ws.Cells.Item(12, 10).Value = 0.0002405259164795302
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 0.2909853333333334
$ws.Cells.Item(12, 14).Value = 0.8729560000000001
$ws.Cells.Item(12, 15).Value = 0.03275300561492853
$ws.Cells.Item(12, 16).Value = 0.03275300561492853
$ws.Cells.Item(12, 17).Value = 0.02507110232977778
$ws.Cells.Item(12, 18).Value = 0.225639920968
$ws.Cells.Item(12, 19).Value = 0.000007877946692989882
$ws.Cells.Item(12, 20).Value = 0.000007877946692989884

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 0.08615933333333332
$ws.Cells.Item(13, 8).Value = 0.258478
$ws.Cells.Item(13, 9).Value = 0.0002405259164795302
$ws.Cells.Item(13, 10).Value = 0.0002405259164795302
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 0.4814053333333333
$ws.Cells.Item(13, 14).Value = 1.444216
$ws.Cells.Item(13, 15).Value = 0.0541864822020464
$ws.Cells.Item(13, 16).Value = 0.05418648220204641
$ws.Cells.Item(13, 17).Value = 0.0414775625831111
$ws.Cells.Item(13, 18).Value = 0.373298063248
$ws.Cells.Item(13, 19).Value = 0.00001303325329244896
$ws.Cells.Item(13, 20).Value = 0.00001303325329244896

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 13.77014866666667
$ws.Cells.Item(14, 8).Value = 41.310446
$ws.Cells.Item(14, 9).Value = 0.03844130983808348
$ws.Cells.Item(14, 10).Value = 0.03844130983808348
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 3.778439
$ws.Cells.Item(14, 14).Value = 11.335317
$ws.Cells.Item(14, 15).Value = 0.4252971528324392
$ws.Cells.Item(14, 16).Value = 0.4252971528324392
$ws.Cells.Item(14, 17).Value = 52.02966675793133
$ws.Cells.Item(14, 18).Value = 468.267000821382
$ws.Cells.Item(14, 19).Value = 0.01634897962528654
$ws.Cells.Item(14, 20).Value = 0.01634897962528654

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 13.77014866666667
$ws.Cells.Item(15, 8).Value = 41.310446
$ws.Cells.Item(15, 9).Value = 0.03844130983808348
$ws.Cells.Item(15, 10).Value = 0.03844130983808348
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 4.333403333333333
$ws.Cells.Item(15, 14).Value = 13.00021
$ws.Cells.Item(15, 15).Value = 0.4877633593505858
$ws.Cells.Item(15, 16).Value = 0.4877633593505858
$ws.Cells.Item(15, 17).Value = 59.67160813262888
$ws.Cells.Item(15, 18).Value = 537.04447319366
$ws.Cells.Item(15, 19).Value = 0.01875026242446032
$ws.Cells.Item(15, 20).Value = 0.01875026242446032

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 13.77014866666667
$ws.Cells.Item(16, 8).Value = 41.310446
$ws.Cells.Item(16, 9).Value = 0.03844130983808348
$ws.Cells.Item(16, 10).Value = 0.03844130983808348
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.2909853333333334
$ws.Cells.Item(16, 14).Value = 0.8729560000000001
$ws.Cells.Item(16, 15).Value = 0.03275300561492853
$ws.Cells.Item(16, 16).Value = 0.03275300561492853
$ws.Cells.Item(16, 17).Value = 4.006911299819556
$ws.Cells.Item(16, 18).Value = 36.062201698376
$ws.Cells.Item(16, 19).Value = 0.001259068436971956
$ws.Cells.Item(16, 20).Value = 0.001259068436971956

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 13.77014866666667
$ws.Cells.Item(17, 8).Value = 41.310446
$ws.Cells.Item(17, 9).Value = 0.03844130983808348
$ws.Cells.Item(17, 10).Value = 0.03844130983808348
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 0.4814053333333333
$ws.Cells.Item(17, 14).Value = 1.444216
$ws.Cells.Item(17, 15).Value = 0.0541864822020464
$ws.Cells.Item(17, 16).Value = 0.05418648220204641
$ws.Cells.Item(17, 17).Value = 6.629023008926222
$ws.Cells.Item(17, 18).Value = 59.661207080336
$ws.Cells.Item(17, 19).Value = 0.002082999351364662
$ws.Cells.Item(17, 20).Value = 0.002082999351364662

Write-Output "Applied NATMI data updates (Dr Hou advice) to 16 rows."